$d = $word.ActiveDocument

$pairs = @(
    @("36÷6=", "89÷4="),
    @("82÷4=", "89÷7="),
    @("92÷8=", "36÷2="),
    @("67÷9=", "28÷5="),
    @("56÷5=", "86÷3="),
    @("48÷5=", "99÷9="),
    @("52÷6=", "38÷5="),
    @("50÷4=", "98÷8="),
    @("93÷3=", "82÷9="),
    @("66÷8=", "80÷5="),
    @("98÷5=", "83÷4="),
    @("97÷6=", "79÷2="),
    @("37÷8=", "18÷3="),
    @("52÷5=", "89÷7="),
    @("76÷9=", "93÷3="),
    @("67÷5=", "64÷9="),
    @("44÷5=", "83÷4="),
    @("61÷8=", "95÷8="),
    @("15÷3=", "81÷3="),
    @("87÷9=", "99÷4="),
    @("31÷2=", "69÷3="),
    @("24÷2=", "59÷6="),
    @("65÷4=", "45÷2="),
    @("86÷4=", "10÷5="),
    @("47÷5=", "27÷6=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $rng = $d.Content
    $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "replacements applied"
